$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes existing rows 5..122 down to 6..123)
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data record
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = "2023-09-21"
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112001
$ws.Range("G5").Value = "Berenjena"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("N5").Value = "`$/caja 60 unidades"
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 167
$ws.Range("Q5").Value = 60
$ws.Range("R5").Value = "Hortaliza"
